$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (113) down through the
# 15 new rows (114-128) so the new cells inherit the same styles
# (date format on col A, 2-decimal number format on col B, default on col C).
$ws.Range("A113:C113").Copy($ws.Range("A114:C128"))

$data = @(
    @(45219.583472222221, 0, 2709.5),
    @(45219.541805555556, 0, 2698),
    @(45219.500138888892, 0, 2684.7),
    @(45219.458472222221, 0, 2662.3),
    @(45219.416805555556, 0, 2571),
    @(45219.375138888892, 0, 2563.4),
    @(45219.333472222221, 0, 2582.6999999999998),
    @(45219.291805555556, 0, 2607.6),
    @(45219.250138888892, 0, 2657.5),
    @(45219.208472222221, 0, 2690.9),
    @(45219.166805555556, 0, 2708.3),
    @(45219.125138888892, 0, 2728.8),
    @(45219.083472222221, 0, 2725.2),
    @(45219.041805555556, 0, 2722),
    @(45219.000138888892, 0, 2721.3)
)

$row = 114
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").Value = $r[2]
    $row++
}

# Restore the view state: scrolled so row 111 is the top visible row,
# with H122 as the active selected cell (as in the saved workbook).
$excel.ActiveWindow.ScrollRow = 111
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("H122").Select()
